$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the bottom of the table: rows 21-29 are removed entirely ---
#     (this also shrinks the sheet's used range / dimension from A1:F29 to A1:F20)
$ws.Rows("21:29").Delete()

# --- 2. Refresh the data for the remaining rows 2-20 ---
# Column A keeps the same 0-based running index (0..18), unchanged by the edit.
# Columns B, C and E get a refreshed list of tickers; D and F are cleared for
# every data row in the new layout.
$colB = @('NSE:BALAJITELE', 'NSE:HUBTOWN', 'NSE:INDIAMART', 'NSE:IVP', 'NSE:M&M', 'NSE:MBAPL', 'NSE:MEDICAMEQ', 'NSE:NUVAMA', 'NSE:PRIVISCL', 'NSE:RML', '', '', '', '', '', '', '', '', '')
$colC = @('NSE:AARTECH', 'NSE:AVTNPL', 'NSE:CTE', 'NSE:DCM', 'NSE:EIDPARRY', 'NSE:GHCL', 'NSE:HATSUN', 'NSE:HDFCGROWTH', 'NSE:IDBI', 'NSE:JISLDVREQS', 'NSE:JSL', 'NSE:KELLTONTEC', 'NSE:MACPOWER', 'NSE:MMP', 'NSE:NITCO', 'NSE:PRSMJOHNSN', 'NSE:PUNJABCHEM', 'NSE:RHL', 'NSE:SAMBHAAV')
$colE = @('NSE:ABB', 'NSE:ADANIENT', 'NSE:CAMS', 'NSE:CONCOR', 'NSE:DIXON', 'NSE:GRASIM', 'NSE:HINDALCO', 'NSE:HUDCO', 'NSE:INDHOTEL', 'NSE:ITC', 'NSE:NBCC', 'NSE:PFC', 'NSE:RECLTD', '', '', '', '', '', '')

for ($i = 0; $i -lt 19; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $colB[$i]
    $ws.Cells.Item($row, 3).Value = $colC[$i]
    $ws.Cells.Item($row, 4).Value = ''
    $ws.Cells.Item($row, 5).Value = $colE[$i]
    $ws.Cells.Item($row, 6).Value = ''
}
